$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner (A1)
$ws.Range("A1").Value = "Datos actualizados a 29 de Abril de 2020 a las 22:22"

# Swap the display order of Dominica / Curazao rows (row 192 <-> 193 country names)
$ws.Range("A192").Value = "Dominica"
$ws.Range("A193").Value = "Curazao"

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 1054166
$ws.Range("C4").Value = 18401
$ws.Range("E4").Value = 848734
$ws.Range("F4").Value = 18665
$ws.Range("G4").Value = 1743
$ws.Range("H4").Value = 61009

# Row 9 (Reino Unido)
$ws.Range("B9").Value = 161173
$ws.Range("C9").Value = 1261
$ws.Range("E9").Value = 34374
$ws.Range("F9").Value = 2415
$ws.Range("G9").Value = 85
$ws.Range("H9").Value = 6399

# Row 18 (Brasil)
$ws.Range("B18").Value = 33062
$ws.Range("C18").Value = 1738
$ws.Range("D18").Value = 8437
$ws.Range("E18").Value = 23546
$ws.Range("G18").Value = 71
$ws.Range("H18").Value = 1079

# Row 25 (Ecuador)
$ws.Range("D25").Value = 13386
$ws.Range("E25").Value = 5677
$ws.Range("F25").Value = 129

# Row 72
$ws.Range("D72").Value = 934
$ws.Range("E72").Value = 811
$ws.Range("G72").Value = 3
$ws.Range("H72").Value = 61

# Row 86
$ws.Range("B86").Value = 1351
$ws.Range("C86").Value = 111
$ws.Range("D86").Value = 313
$ws.Range("E86").Value = 1031

# Row 91
$ws.Range("B91").Value = 980
$ws.Range("C91").Value = 5
$ws.Range("D91").Value = 294
$ws.Range("E91").Value = 646
$ws.Range("F91").Value = 20

# Row 127
$ws.Range("D127").Value = 258
$ws.Range("E127").Value = 34
$ws.Range("F127").Value = 21

# Row 192 (province data values)
$ws.Range("D192").Value = 13
$ws.Range("E192").Value = 3
$ws.Range("H192").Value = 0

# Row 193 (province data values)
$ws.Range("E193").Value = 2
$ws.Range("H193").Value = 1
